# Add files via upload / Translated to English
# Translate the Chinese assembly notes (rows 7, 35-40 of column A) to English,
# matching the published BOM update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: short PCB-thickness note -> gets a leading blank line and becomes a
# wrapped, taller row in the English revision.
$ws.Range("A7").Value = "`nPCB thickness 1mm"
$ws.Range("A7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 30

# Rows 35-40: translate the jumper / mode explanation notes.
$ws.Range("A35").Value = "RI-R6 can be used without using it (if used, it can increase the compatibility of some TF cards)"
$ws.Range("A36").Value = "R8-R11 is only used to increase the thickness to fit the SD card slot of RHEA/PHOEBE. If you are making your own for fenrir only, you don’t need to use it."
$ws.Range("A37").Value = "MODE jumper: This jumper sets the function of SW2"
$ws.Range("A38").Value = "BUTTON mode is suitable for RHEA/PHOEBE. In this case, the FN pad needs to be connected to the button switch of RHEA/PHOEBE."
$ws.Range("A39").Value = "The DOOR mode simulates the opening and closing of the CD compartment cover, which is applicable to some FENRIRs. However, since FENRIR only retains the hardware, the software currently does not enable this function, so it will not work now. The author CED needs to add this feature."
$ws.Range("A40").Value = "The factory settings of the merchant can be defined in DOOR mode, and RHEA/PHOEBE is compatible with this mode"

# Matches the saved file's final selection (scrolled further down the sheet).
[void]$ws.Range("C40").Select()
